# Apply Priority / Ready to be planned / Planned Sprint / Status values
# to the "Product Backlog Template" sheet rows that already contain a
# user story (columns D, E, F, G).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Product Backlog Template")

# row -> (Priority, Ready to be planned, Planned Sprint, Status)
$rows = @{
    2  = @("High",   "Yes", 1, "In progress")
    3  = @("High",   "Yes", 1, "In progress")
    4  = @("High",   "Yes", 1, "In progress")
    6  = @("High",   "Yes", 1, "Open")
    7  = @("Medium",  "No", 2, "Todo")
    8  = @("Low",     "No", 3, "Todo")
    10 = @("Medium",  "No", 3, "Todo")
    11 = @("Medium",  "No", 3, "Todo")
    13 = @("High",   "Yes", 1, "Open")
    14 = @("High",   "Yes", 1, "Open")
    17 = @("Medium", "Yes", 2, "Todo")
    18 = @("Medium", "Yes", 2, "Todo")
    19 = @("Medium", "Yes", 2, "Todo")
    20 = @("Medium", "Yes", 2, "Todo")
    23 = @("High",   "Yes", 1, "Open")
    24 = @("High",   "Yes", 1, "Open")
    25 = @("High",   "Yes", 1, "Open")
    27 = @("Medium",  "No", 3, "Todo")
    28 = @("Medium",  "No", 3, "Todo")
    30 = @("Low",     "No", 3, "Todo")
    31 = @("Low",     "No", 3, "Todo")
}

foreach ($r in $rows.Keys) {
    $vals = $rows[$r]
    $ws.Range("D$r").Value = $vals[0]
    $ws.Range("E$r").Value = $vals[1]
    $ws.Range("F$r").Value = $vals[2]
    $ws.Range("G$r").Value = $vals[3]
    $ws.Range("D$r:G$r").HorizontalAlignment = -4108
}
